# Update the 5x3 lattice-multiplication exercise table: each cell
# gets a new "A x B" problem with matching partial-product digits and
# lattice row labels. Cell XML structure (single run: header, factor
# line, dashed rule, two lattice-row labels, separated by <w:br/>) is
# preserved; only the text content changes.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$cell = $t.Cell(1, 1)
[void]$cell.Range.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r><w:rPr><w:sz w:val='32'/></w:rPr><w:t>89 x 99</w:t><w:br/><w:t xml:space='preserve'>  9    9</w:t><w:br/><w:t xml:space='preserve'>  ----</w:t><w:br/><w:t>8|    |</w:t><w:br/><w:t>9|    |</w:t></w:r></w:p>")

$cell = $t.Cell(1, 2)
[void]$cell.Range.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r><w:rPr><w:sz w:val='32'/></w:rPr><w:t>45 x 60</w:t><w:br/><w:t xml:space='preserve'>  6    0</w:t><w:br/><w:t xml:space='preserve'>  ----</w:t><w:br/><w:t>4|    |</w:t><w:br/><w:t>5|    |</w:t></w:r></w:p>")

$cell = $t.Cell(1, 3)
[void]$cell.Range.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r><w:rPr><w:sz w:val='32'/></w:rPr><w:t>43 x 31</w:t><w:br/><w:t xml:space='preserve'>  3    1</w:t><w:br/><w:t xml:space='preserve'>  ----</w:t><w:br/><w:t>4|    |</w:t><w:br/><w:t>3|    |</w:t></w:r></w:p>")

$cell = $t.Cell(2, 1)
[void]$cell.Range.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r><w:rPr><w:sz w:val='32'/></w:rPr><w:t>61 x 47</w:t><w:br/><w:t xml:space='preserve'>  4    7</w:t><w:br/><w:t xml:space='preserve'>  ----</w:t><w:br/><w:t>6|    |</w:t><w:br/><w:t>1|    |</w:t></w:r></w:p>")

$cell = $t.Cell(2, 2)
[void]$cell.Range.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r><w:rPr><w:sz w:val='32'/></w:rPr><w:t>29 x 10</w:t><w:br/><w:t xml:space='preserve'>  1    0</w:t><w:br/><w:t xml:space='preserve'>  ----</w:t><w:br/><w:t>2|    |</w:t><w:br/><w:t>9|    |</w:t></w:r></w:p>")

$cell = $t.Cell(2, 3)
[void]$cell.Range.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r><w:rPr><w:sz w:val='32'/></w:rPr><w:t>63 x 38</w:t><w:br/><w:t xml:space='preserve'>  3    8</w:t><w:br/><w:t xml:space='preserve'>  ----</w:t><w:br/><w:t>6|    |</w:t><w:br/><w:t>3|    |</w:t></w:r></w:p>")

$cell = $t.Cell(3, 1)
[void]$cell.Range.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r><w:rPr><w:sz w:val='32'/></w:rPr><w:t>74 x 16</w:t><w:br/><w:t xml:space='preserve'>  1    6</w:t><w:br/><w:t xml:space='preserve'>  ----</w:t><w:br/><w:t>7|    |</w:t><w:br/><w:t>4|    |</w:t></w:r></w:p>")

$cell = $t.Cell(3, 2)
[void]$cell.Range.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r><w:rPr><w:sz w:val='32'/></w:rPr><w:t>53 x 30</w:t><w:br/><w:t xml:space='preserve'>  3    0</w:t><w:br/><w:t xml:space='preserve'>  ----</w:t><w:br/><w:t>5|    |</w:t><w:br/><w:t>3|    |</w:t></w:r></w:p>")

$cell = $t.Cell(3, 3)
[void]$cell.Range.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r><w:rPr><w:sz w:val='32'/></w:rPr><w:t>15 x 71</w:t><w:br/><w:t xml:space='preserve'>  7    1</w:t><w:br/><w:t xml:space='preserve'>  ----</w:t><w:br/><w:t>1|    |</w:t><w:br/><w:t>5|    |</w:t></w:r></w:p>")

$cell = $t.Cell(4, 1)
[void]$cell.Range.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r><w:rPr><w:sz w:val='32'/></w:rPr><w:t>60 x 20</w:t><w:br/><w:t xml:space='preserve'>  2    0</w:t><w:br/><w:t xml:space='preserve'>  ----</w:t><w:br/><w:t>6|    |</w:t><w:br/><w:t>0|    |</w:t></w:r></w:p>")

$cell = $t.Cell(4, 2)
[void]$cell.Range.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r><w:rPr><w:sz w:val='32'/></w:rPr><w:t>55 x 65</w:t><w:br/><w:t xml:space='preserve'>  6    5</w:t><w:br/><w:t xml:space='preserve'>  ----</w:t><w:br/><w:t>5|    |</w:t><w:br/><w:t>5|    |</w:t></w:r></w:p>")

$cell = $t.Cell(4, 3)
[void]$cell.Range.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r><w:rPr><w:sz w:val='32'/></w:rPr><w:t>66 x 80</w:t><w:br/><w:t xml:space='preserve'>  8    0</w:t><w:br/><w:t xml:space='preserve'>  ----</w:t><w:br/><w:t>6|    |</w:t><w:br/><w:t>6|    |</w:t></w:r></w:p>")

$cell = $t.Cell(5, 1)
[void]$cell.Range.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r><w:rPr><w:sz w:val='32'/></w:rPr><w:t>94 x 81</w:t><w:br/><w:t xml:space='preserve'>  8    1</w:t><w:br/><w:t xml:space='preserve'>  ----</w:t><w:br/><w:t>9|    |</w:t><w:br/><w:t>4|    |</w:t></w:r></w:p>")

$cell = $t.Cell(5, 2)
[void]$cell.Range.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r><w:rPr><w:sz w:val='32'/></w:rPr><w:t>47 x 53</w:t><w:br/><w:t xml:space='preserve'>  5    3</w:t><w:br/><w:t xml:space='preserve'>  ----</w:t><w:br/><w:t>4|    |</w:t><w:br/><w:t>7|    |</w:t></w:r></w:p>")

$cell = $t.Cell(5, 3)
[void]$cell.Range.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r><w:rPr><w:sz w:val='32'/></w:rPr><w:t>46 x 72</w:t><w:br/><w:t xml:space='preserve'>  7    2</w:t><w:br/><w:t xml:space='preserve'>  ----</w:t><w:br/><w:t>4|    |</w:t><w:br/><w:t>6|    |</w:t></w:r></w:p>")
